# Uppercase the month header labels in the calendar worksheet.
# Each month name appears twice per month block: once in column E and
# once in column M, on the same header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$months = @{
    2  = "JANUARY"
    10 = "FEBRUARY"
    18 = "MARCH"
    26 = "APRIL"
    34 = "MAY"
    42 = "JUNE"
    50 = "JULY"
    58 = "AUGUST"
    66 = "SEPTEMBER"
    75 = "OCTOBER"
    83 = "NOVEMBER"
    91 = "DECEMBER"
}

foreach ($row in $months.Keys) {
    $name = $months[$row]
    $ws.Range("E$row").Value = $name
    $ws.Range("M$row").Value = $name
}
